# reset position for all scene
# Updates the RelivePos ("复活坐标列表") column (G) for rows 11-16 on the
# Scene sheet with new per-scene revive coordinates, and restores the
# sheet's view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G11").Value = "89,104,0"
$ws.Range("G12").Value = "89,102,0"
$ws.Range("G13").Value = "89,104,0"
$ws.Range("G14").Value = "100,115,0"
$ws.Range("G15").Value = "95,102,0"
$ws.Range("G16").Value = "88,105,0"
